$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.777.03"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.542.34"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.78"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.10"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.41"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.30"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.934.63"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.42"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.552.01"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.816"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.779.91"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0951"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.83"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.05"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.88"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.75"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.13"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.27"
$ws.Range("E30").Value = "  -4.45%  "
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.80"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.73"
$ws.Range("E33").Value = "  +4.74%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.67"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0793"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.87"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.76"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.999.69"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.787.24"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.191"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.78"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.36"
